$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    11  = 6.58
    12  = 4.951000000000001
    15  = 5.135
    27  = 5.424
    28  = 5.81
    31  = 5.779
    32  = 6.589
    36  = 8.443
    38  = 5.511000000000001
    46  = 6.706999999999999
    54  = 5.176
    55  = 4.722
    56  = 4.654000000000001
    67  = 5.516
    69  = 5.366999999999999
    72  = 5.845000000000001
    73  = 6.935
    83  = 5.827
    86  = 5.02
    91  = 5.667
    93  = 5.423999999999999
    99  = 5.699000000000001
    104 = 8.093999999999999
    105 = 8.076000000000001
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 2).Value = $updates[$row]
}
